$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$columns = @("C", "E", "F", "G", "H", "I", "J", "K", "L", "M", "N", "O", "P", "Q", "R", "S", "T", "U", "V", "W", "X", "Y")

$data = @(
    @(2, 0.4439959526062012, 67.76515235949591, 0.002251725967537097, 0.001994023637924339, 0.001770823232218035, 0.001696340112274568, 0.001696340112274568, 0.001652166767338213, 0.001640092300817573, 0.001569885235033061, 0.001522893761134875, 0.001498910872656215, 0.001494442772140172, 0.001480412947147515, 0.001435774385189924, 0.001406529534540266, 0.001375722688603682, 0.001369069212867554, 0.001350552546095126, 0.001320958135662688, 0.001320958135662688, 0.001320958135662688),
    @(3, 0.3639762401580811, 67.7677627982157, 0.002251725967537097, 0.001905245888792103, 0.001796413565771643, 0.001724464479688052, 0.001680122351024738, 0.001605002834890045, 0.001552212583888019, 0.00145065100351908, 0.001432332768449217, 0.001432332768449217, 0.001409792264924772, 0.001373399205807708, 0.001331799429094263, 0.001331799429094263, 0.001331799429094263, 0.001325913478444472, 0.001325913478444472, 0.001323582746114984, 0.001323582746114984, 0.001321009021407713),
    @(4, 0.3619959354400635, 68.07349966671245, 0.002251725967537097, 0.001882451125012616, 0.001703665685309692, 0.001703665685309692, 0.001654810202386467, 0.00154177382804715, 0.00154177382804715, 0.001437088300089596, 0.001437088300089596, 0.001432925978414713, 0.001432925978414713, 0.001432925978414713, 0.001401134136565383, 0.001401134136565383, 0.001367028864902292, 0.001347986721313227, 0.001347986721313227, 0.001342204894243936, 0.001327514534611951, 0.001326968804419346),
    @(5, 0.3850030899047852, 68.62686388740258, 0.002168751758780943, 0.001848732562155039, 0.001795491960731428, 0.001732684028484852, 0.001684583754126931, 0.001684583754126931, 0.001578225675324287, 0.001551009421582137, 0.00151743919092832, 0.001482640172362649, 0.001479093105384426, 0.0014303195172292, 0.001407189344119073, 0.001407189344119073, 0.001407189344119073, 0.00139213927090051, 0.001366885030886754, 0.001361060515916043, 0.001349760237456449, 0.001337755631333383),
    @(6, 0.3379988670349121, 66.80591782280317, 0.002204778545087303, 0.00186528695705306, 0.001672953348297419, 0.001600387213193439, 0.001600387213193439, 0.00158363099104048, 0.001458513245430779, 0.001458513245430779, 0.001425815106705951, 0.001419779103914735, 0.001418726849932279, 0.001418726849932279, 0.001415889140880147, 0.001388645477659506, 0.001361948391625874, 0.001353016288948954, 0.00133478762765907, 0.001318786873233438, 0.001305618593374235, 0.001302259606682322),
    @(7, 0.5020124912261963, 69.38432292406105, 0.002225726544130102, 0.001887817465543896, 0.001721600267429329, 0.001620301991031341, 0.001579629437243875, 0.001572164832052354, 0.001480817476018366, 0.001480817476018366, 0.001449421184235128, 0.001443042738088163, 0.001428876708005882, 0.001386100609925477, 0.001386100609925477, 0.001369916374965719, 0.001369916374965719, 0.001361923397848849, 0.001359403267458448, 0.001357159529092399, 0.001354799215832774, 0.001352520914699046),
    @(8, 0.4069933891296387, 69.05400904423004, 0.002251725967537097, 0.001808058442435551, 0.001684724324544799, 0.001684724324544799, 0.001601466755486608, 0.001572679556720189, 0.001549519067305342, 0.001489991326940234, 0.001489991326940234, 0.001439448748918395, 0.001432058428655797, 0.001432058428655797, 0.001430151461905308, 0.001425147593333897, 0.00142097469155842, 0.001380355537067801, 0.001376123048518429, 0.001364499286267277, 0.001359155639952999, 0.001346082047645809),
    @(9, 0.3819990158081055, 72.4433696700853, 0.002233534034567199, 0.001913506730443277, 0.001825601836752712, 0.00172829899103374, 0.001706608045708165, 0.001665378163008859, 0.001646406816097274, 0.001646406816097274, 0.001518935396185578, 0.001518935396185578, 0.001518935396185578, 0.001473765359558874, 0.001473765359558874, 0.001445096344670592, 0.001445096344670592, 0.001445096344670592, 0.001426308546948032, 0.001426308546948032, 0.001420472305206764, 0.001412151455557218),
    @(10, 0.3600194454193115, 64.49166949392566, 0.002251725967537097, 0.001909053467736253, 0.00176948621282511, 0.001604848932859711, 0.001558804230442203, 0.001441795458406232, 0.001441795458406232, 0.001370923270996674, 0.001354413872092431, 0.001354413872092431, 0.001354413872092431, 0.001354413872092431, 0.00130281713131175, 0.00130281713131175, 0.001300349323208787, 0.001279586579626571, 0.001279586579626571, 0.001262349777819599, 0.001262349777819599, 0.001257147553487829),
    @(11, 0.4619977474212646, 68.85805644363791, 0.002230626523665372, 0.001848276359618658, 0.001750678026349059, 0.001746339438997173, 0.001698384289887465, 0.001566395556625993, 0.001533066924512063, 0.001533066924512063, 0.001468617783662127, 0.001371934846889101, 0.00136448471432072, 0.00136448471432072, 0.00136448471432072, 0.00136448471432072, 0.001362897150696523, 0.001359522165112352, 0.001347864070805662, 0.001346556143593116, 0.001343654981714986, 0.001342262308842844),
)

foreach ($row in $data) {
    $r = $row[0]
    for ($i = 0; $i -lt $columns.Length; $i++) {
        $col = $columns[$i]
        $val = $row[$i + 1]
        $ws.Range("$col$r").Value = $val
    }
}
